$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "hey I made a change"
$ws.Range("A2").Select()
